# Auto-generated edit script applying updated cosinor analysis results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"24.40000000000038"
$ws.Range("H2").Value = [double]"3.983480212355062e-13"
$ws.Range("I2").Value = [double]"3.983480212355062e-13"
$ws.Range("L2").Value = [double]"56.91859277802148"
$ws.Range("M2").Value = "[45.58468841485261, 68.25249714119036]"
$ws.Range("N2").Value = [double]"3.632649736573512e-13"
$ws.Range("O2").Value = [double]"3.632649736573512e-13"
$ws.Range("P2").Value = [double]"1.578658170272348"
$ws.Range("Q2").Value = "[1.3522370781217319, 1.8050792624229643]"
$ws.Range("T2").Value = [double]"56.03680821834028"
$ws.Range("U2").Value = "[48.421572637350444, 63.65204379933011]"
$ws.Range("X2").Value = [double]"18.26946946946975"
$ws.Range("Y2").Value = [double]"17.39019019019046"
$ws.Range("Z2").Value = [double]"19.14874874874905"
$ws.Range("F3").Value = [double]"24.40000000000038"
$ws.Range("H3").Value = [double]"6.828759779864413e-12"
$ws.Range("I3").Value = [double]"6.828759779864413e-12"
$ws.Range("L3").Value = [double]"61.0169808986368"
$ws.Range("M3").Value = "[47.44236313590024, 74.59159866137337]"
$ws.Range("N3").Value = [double]"1.06581410364015e-11"
$ws.Range("O3").Value = [double]"1.06581410364015e-11"
$ws.Range("P3").Value = [double]"1.528342372016656"
$ws.Range("Q3").Value = "[1.2641844311742716, 1.7925003128590404]"
$ws.Range("R3").Value = [double]"3.552713678800501e-15"
$ws.Range("S3").Value = [double]"3.552713678800501e-15"
$ws.Range("T3").Value = [double]"52.01896668952315"
$ws.Range("U3").Value = "[43.083505725173126, 60.95442765387318]"
$ws.Range("V3").Value = [double]"2.886579864025407e-15"
$ws.Range("W3").Value = [double]"2.886579864025407e-15"
$ws.Range("X3").Value = [double]"18.46486486486515"
$ws.Range("Y3").Value = [double]"17.43903903903931"
$ws.Range("Z3").Value = [double]"19.49069069069099"
$ws.Range("F4").Value = [double]"24.40000000000038"
$ws.Range("H4").Value = [double]"2.433031554005538e-11"
$ws.Range("I4").Value = [double]"2.433031554005538e-11"
$ws.Range("L4").Value = [double]"59.86629657083881"
$ws.Range("M4").Value = "[45.7506867019806, 73.98190643969701]"
$ws.Range("N4").Value = [double]"5.679900993982301e-11"
$ws.Range("O4").Value = [double]"5.679900993982301e-11"
$ws.Range("P4").Value = [double]"1.515763422452732"
$ws.Range("Q4").Value = "[1.2390265320464247, 1.7925003128590395]"
$ws.Range("R4").Value = [double]"2.19824158875781e-14"
$ws.Range("S4").Value = [double]"2.19824158875781e-14"
$ws.Range("T4").Value = [double]"56.62638572587814"
$ws.Range("U4").Value = "[47.48263348402272, 65.77013796773358]"
$ws.Range("V4").Value = [double]"4.440892098500626e-16"
$ws.Range("W4").Value = [double]"4.440892098500626e-16"
$ws.Range("X4").Value = [double]"18.513713713714"
$ws.Range("Y4").Value = [double]"17.43903903903932"
$ws.Range("Z4").Value = [double]"19.58838838838869"
$ws.Range("F5").Value = [double]"24.40000000000038"
$ws.Range("H5").Value = [double]"1.319833131674386e-12"
$ws.Range("I5").Value = [double]"1.319833131674386e-12"
$ws.Range("L5").Value = [double]"57.22955806171375"
$ws.Range("M5").Value = "[43.70414940925916, 70.75496671416835]"
$ws.Range("N5").Value = [double]"6.065525859355603e-11"
$ws.Range("O5").Value = [double]"6.065525859355603e-11"
$ws.Range("P5").Value = [double]"1.264184431174272"
$ws.Range("Q5").Value = "[1.0126054398958093, 1.5157634224527339]"
$ws.Range("R5").Value = [double]"3.563815909046752e-13"
$ws.Range("S5").Value = [double]"3.563815909046752e-13"
$ws.Range("T5").Value = [double]"56.6572370284004"
$ws.Range("U5").Value = "[48.734331685185026, 64.58014237161578]"
$ws.Range("V5").Value = [double]"0"
$ws.Range("W5").Value = [double]"0"
$ws.Range("X5").Value = [double]"19.49069069069099"
$ws.Range("Y5").Value = [double]"18.513713713714"
$ws.Range("Z5").Value = [double]"20.46766766766798"
$ws.Range("F6").Value = [double]"24.40000000000038"
$ws.Range("H6").Value = [double]"1.570810148621149e-11"
$ws.Range("I6").Value = [double]"1.570810148621149e-11"
$ws.Range("L6").Value = [double]"55.50618677837979"
$ws.Range("M6").Value = "[41.16408316910234, 69.84829038765723]"
$ws.Range("N6").Value = [double]"6.865645829634559e-10"
$ws.Range("O6").Value = [double]"6.865645829634559e-10"
$ws.Range("P6").Value = [double]"1.213868632918578"
$ws.Range("Q6").Value = "[0.9245527929483472, 1.5031844728888082]"
$ws.Range("R6").Value = [double]"7.688494285673642e-11"
$ws.Range("S6").Value = [double]"7.688494285673642e-11"
$ws.Range("T6").Value = [double]"48.08665347446178"
$ws.Range("U6").Value = "[39.76342581754097, 56.4098811313826]"
$ws.Range("V6").Value = [double]"3.552713678800501e-15"
$ws.Range("W6").Value = [double]"3.552713678800501e-15"
$ws.Range("X6").Value = [double]"19.68608608608639"
$ws.Range("Y6").Value = [double]"18.56256256256285"
$ws.Range("Z6").Value = [double]"20.80960960960993"
$ws.Range("F7").Value = [double]"24.40000000000038"
$ws.Range("H7").Value = [double]"7.105427357601002e-15"
$ws.Range("I7").Value = [double]"7.105427357601002e-15"
$ws.Range("L7").Value = [double]"56.70771118068951"
$ws.Range("M7").Value = "[44.38697194116, 69.02845042021903]"
$ws.Range("N7").Value = [double]"5.28577182024037e-12"
$ws.Range("O7").Value = [double]"5.28577182024037e-12"
$ws.Range("P7").Value = [double]"1.125815985971117"
$ws.Range("Q7").Value = "[0.8993948938205012, 1.3522370781217319]"
$ws.Range("R7").Value = [double]"4.96713781217295e-13"
$ws.Range("S7").Value = [double]"4.96713781217295e-13"
$ws.Range("T7").Value = [double]"53.12572562974447"
$ws.Range("U7").Value = "[46.43383443543784, 59.81761682405109]"
$ws.Range("X7").Value = [double]"20.02802802802834"
$ws.Range("Y7").Value = [double]"19.14874874874905"
$ws.Range("Z7").Value = [double]"20.90730730730763"
$ws.Range("F8").Value = [double]"23.84000000000029"
$ws.Range("H8").Value = [double]"3.693934047532821e-12"
$ws.Range("I8").Value = [double]"3.693934047532821e-12"
$ws.Range("L8").Value = [double]"60.13672467048717"
$ws.Range("M8").Value = "[43.15918609523747, 77.11426324573686]"
$ws.Range("N8").Value = [double]"6.446097433610021e-09"
$ws.Range("O8").Value = [double]"6.446097433610021e-09"
$ws.Range("P8").Value = [double]"0.823921196436963"
$ws.Range("Q8").Value = "[0.5723422051585008, 1.0755001877154253]"
$ws.Range("R8").Value = [double]"4.055176083817003e-08"
$ws.Range("S8").Value = [double]"4.055176083817003e-08"
$ws.Range("T8").Value = [double]"53.73423469926428"
$ws.Range("U8").Value = "[45.16851649871967, 62.2999528998089]"
$ws.Range("V8").Value = [double]"2.220446049250313e-16"
$ws.Range("W8").Value = [double]"2.220446049250313e-16"
$ws.Range("X8").Value = [double]"20.71383383383408"
$ws.Range("Y8").Value = [double]"19.75927927927952"
$ws.Range("Z8").Value = [double]"21.66838838838865"
$ws.Range("F9").Value = [double]"23.84000000000029"
$ws.Range("H9").Value = [double]"2.752020833440838e-12"
$ws.Range("I9").Value = [double]"2.752020833440838e-12"
$ws.Range("L9").Value = [double]"50.90265667506133"
$ws.Range("M9").Value = "[37.27136251503491, 64.53395083508775]"
$ws.Range("N9").Value = [double]"1.730604104466238e-09"
$ws.Range("O9").Value = [double]"1.730604104466238e-09"
$ws.Range("P9").Value = [double]"0.7987632973091161"
$ws.Range("Q9").Value = "[0.5220264069028078, 1.0755001877154244]"
$ws.Range("R9").Value = [double]"5.903508242255384e-07"
$ws.Range("S9").Value = [double]"5.903508242255384e-07"
$ws.Range("T9").Value = [double]"49.25143441117856"
$ws.Range("U9").Value = "[42.07294031041135, 56.42992851194577]"
$ws.Range("V9").Value = [double]"0"
$ws.Range("W9").Value = [double]"0"
$ws.Range("X9").Value = [double]"20.80928928928954"
$ws.Range("Y9").Value = [double]"19.75927927927952"
$ws.Range("Z9").Value = [double]"21.85929929929956"
$ws.Range("F10").Value = [double]"23.84000000000029"
$ws.Range("H10").Value = [double]"2.343680804983705e-12"
$ws.Range("I10").Value = [double]"2.343680804983705e-12"
$ws.Range("L10").Value = [double]"57.08392722972398"
$ws.Range("M10").Value = "[41.98882842236772, 72.17902603708023]"
$ws.Range("N10").Value = [double]"1.253106063359155e-09"
$ws.Range("O10").Value = [double]"1.253106063359155e-09"
$ws.Range("P10").Value = [double]"1.025184389459731"
$ws.Range("Q10").Value = "[0.748447499053424, 1.301921279866039]"
$ws.Range("R10").Value = [double]"2.119589392890475e-09"
$ws.Range("S10").Value = [double]"2.119589392890475e-09"
$ws.Range("T10").Value = [double]"53.8701049906449"
$ws.Range("U10").Value = "[45.87188798588454, 61.86832199540525]"
$ws.Range("V10").Value = [double]"0"
$ws.Range("W10").Value = [double]"0"
$ws.Range("X10").Value = [double]"19.95019019019043"
$ws.Range("Y10").Value = [double]"18.90018018018041"
$ws.Range("Z10").Value = [double]"21.00020020020045"
$ws.Range("F11").Value = [double]"23.84000000000029"
$ws.Range("H11").Value = [double]"6.6472383153382e-12"
$ws.Range("I11").Value = [double]"6.6472383153382e-12"
$ws.Range("L11").Value = [double]"58.43478429974041"
$ws.Range("M11").Value = "[42.37344959212358, 74.49611900735725]"
$ws.Range("N11").Value = [double]"3.335702603024515e-09"
$ws.Range("O11").Value = [double]"3.335702603024515e-09"
$ws.Range("P11").Value = [double]"0.9874475407679633"
$ws.Range("Q11").Value = "[0.698131700797731, 1.2767633807381955]"
$ws.Range("R11").Value = [double]"1.566022778121123e-08"
$ws.Range("S11").Value = [double]"1.566022778121123e-08"
$ws.Range("T11").Value = [double]"56.31280654989749"
$ws.Range("U11").Value = "[47.838826652203196, 64.7867864475918]"
$ws.Range("X11").Value = [double]"20.09337337337362"
$ws.Range("Y11").Value = [double]"18.99563563563586"
$ws.Range("Z11").Value = [double]"21.19111111111137"
$ws.Range("F12").Value = [double]"23.84000000000029"
$ws.Range("H12").Value = [double]"1.857042297714884e-10"
$ws.Range("I12").Value = [double]"1.857042297714884e-10"
$ws.Range("L12").Value = [double]"56.26438215785407"
$ws.Range("M12").Value = "[38.432490904184114, 74.09627341152401]"
$ws.Range("N12").Value = [double]"9.263108324830682e-08"
$ws.Range("O12").Value = [double]"9.263108324830682e-08"
$ws.Range("P12").Value = [double]"0.9119738433844242"
$ws.Range("Q12").Value = "[0.5975001042863459, 1.2264475824825025]"
$ws.Range("R12").Value = [double]"5.375660059314669e-07"
$ws.Range("S12").Value = [double]"5.375660059314669e-07"
$ws.Range("T12").Value = [double]"56.98515265234647"
$ws.Range("U12").Value = "[47.85352558627933, 66.11677971841361]"
$ws.Range("V12").Value = [double]"2.220446049250313e-16"
$ws.Range("W12").Value = [double]"2.220446049250313e-16"
$ws.Range("X12").Value = [double]"20.37973973973999"
$ws.Range("Y12").Value = [double]"19.18654654654678"
$ws.Range("Z12").Value = [double]"21.5729329329332"

Write-Output "Applied all cell updates."
